$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-9 with the new TPM-derived values
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nlgn2"
$ws.Range("C2").Value = "Nrxn1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.049772999999999
$ws.Range("H2").Value = 15.149319
$ws.Range("I2").Value = 0.1400646900514762
$ws.Range("J2").Value = 0.1400646900514762
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4427803333333333
$ws.Range("N2").Value = 1.328341
$ws.Range("O2").Value = 0.8276247286611124
$ws.Range("P2").Value = 0.8276247286611124
$ws.Range("Q2").Value = 2.235940172197666
$ws.Range("R2").Value = 20.123461549779
$ws.Range("S2").Value = 0.1159210010988558
$ws.Range("T2").Value = 0.1159210010988558
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nlgn2"
$ws.Range("C3").Value = "Nrxn1"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.049772999999999
$ws.Range("H3").Value = 15.149319
$ws.Range("I3").Value = 0.1400646900514762
$ws.Range("J3").Value = 0.1400646900514762
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.092221
$ws.Range("N3").Value = 0.276663
$ws.Range("O3").Value = 0.1723752713388876
$ws.Range("P3").Value = 0.1723752713388876
$ws.Range("Q3").Value = 0.4656951158329999
$ws.Range("R3").Value = 4.191256042497
$ws.Range("S3").Value = 0.02414368895262041
$ws.Range("T3").Value = 0.0241436889526204
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Nlgn2"
$ws.Range("C4").Value = "Nrxn1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 20.420946
$ws.Range("H4").Value = 61.262838
$ws.Range("I4").Value = 0.5664122866607931
$ws.Range("J4").Value = 0.5664122866607931
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4427803333333333
$ws.Range("N4").Value = 1.328341
$ws.Range("O4").Value = 0.8276247286611124
$ws.Range("P4").Value = 0.8276247286611124
$ws.Range("Q4").Value = 9.041993276862
$ws.Range("R4").Value = 81.377939491758
$ws.Range("S4").Value = 0.4687768150579591
$ws.Range("T4").Value = 0.4687768150579591
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Nlgn2"
$ws.Range("C5").Value = "Nrxn1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 20.420946
$ws.Range("H5").Value = 61.262838
$ws.Range("I5").Value = 0.5664122866607931
$ws.Range("J5").Value = 0.5664122866607931
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.092221
$ws.Range("N5").Value = 0.276663
$ws.Range("O5").Value = 0.1723752713388876
$ws.Range("P5").Value = 0.1723752713388876
$ws.Range("Q5").Value = 1.883240061066
$ws.Range("R5").Value = 16.949160549594
$ws.Range("S5").Value = 0.09763547160283402
$ws.Range("T5").Value = 0.09763547160283401
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Nlgn2"
$ws.Range("C6").Value = "Nrxn1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 10.258872
$ws.Range("H6").Value = 30.776616
$ws.Range("I6").Value = 0.2845485781158417
$ws.Range("J6").Value = 0.2845485781158416
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4427803333333333
$ws.Range("N6").Value = 1.328341
$ws.Range("O6").Value = 0.8276247286611124
$ws.Range("P6").Value = 0.8276247286611124
$ws.Range("Q6").Value = 4.542426763784
$ws.Range("R6").Value = 40.881840874056
$ws.Range("S6").Value = 0.2354994397540288
$ws.Range("T6").Value = 0.2354994397540288
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Nlgn2"
$ws.Range("C7").Value = "Nrxn1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.258872
$ws.Range("H7").Value = 30.776616
$ws.Range("I7").Value = 0.2845485781158417
$ws.Range("J7").Value = 0.2845485781158416
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.092221
$ws.Range("N7").Value = 0.276663
$ws.Range("O7").Value = 0.1723752713388876
$ws.Range("P7").Value = 0.1723752713388876
$ws.Range("Q7").Value = 0.946083434712
$ws.Range("R7").Value = 8.514750912407999
$ws.Range("S7").Value = 0.04904913836181287
$ws.Range("T7").Value = 0.04904913836181286
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Nlgn2"
$ws.Range("C8").Value = "Nrxn1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.323557
$ws.Range("H8").Value = 0.9706710000000001
$ws.Range("I8").Value = 0.008974445171889013
$ws.Range("J8").Value = 0.008974445171889013
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4427803333333333
$ws.Range("N8").Value = 1.328341
$ws.Range("O8").Value = 0.8276247286611124
$ws.Range("P8").Value = 0.8276247286611124
$ws.Range("Q8").Value = 0.1432646763123334
$ws.Range("R8").Value = 1.289382086811
$ws.Range("S8").Value = 0.007427472750268674
$ws.Range("T8").Value = 0.007427472750268674
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Nlgn2"
$ws.Range("C9").Value = "Nrxn1"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.323557
$ws.Range("H9").Value = 0.9706710000000001
$ws.Range("I9").Value = 0.008974445171889013
$ws.Range("J9").Value = 0.008974445171889013
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.092221
$ws.Range("N9").Value = 0.276663
$ws.Range("O9").Value = 0.1723752713388876
$ws.Range("P9").Value = 0.1723752713388876
$ws.Range("Q9").Value = 0.029838750097
$ws.Range("R9").Value = 0.268548750873
$ws.Range("S9").Value = 0.001546972421620339
$ws.Range("T9").Value = 0.001546972421620338

# Remove now-obsolete rows 10-13 (data reduced from 12 to 8 rows)
$ws.Rows("10:13").Delete()

Write-Host "Done. UsedRange: $($ws.UsedRange.Address())"
